# Update cryptos list to the latest scraped snapshot.
# Price (column D) is stored as text in the source data (it mirrors
# the raw scraped formatting, e.g. "65.816.01" or "0.0₃0748"), so a
# leading apostrophe forces Excel to keep numeric-looking prices
# ("1.00", "564.68", ...) as text instead of coercing them to numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''65.816.01'
$ws.Range('E2').Value = '  +0.78%  '
$ws.Range('D3').Value = '''3.389.83'
$ws.Range('E3').Value = '  -0.43%  '
$ws.Range('D4').Value = '''1.00'
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').Value = '''564.68'
$ws.Range('E5').Value = '  +0.60%  '
$ws.Range('D6').Value = '''176.76'
$ws.Range('E6').Value = '  +0.34%  '
$ws.Range('D7').Value = '''0.630'
$ws.Range('E7').Value = '  +0.52%  '
$ws.Range('D8').Value = '''3.383.90'
$ws.Range('E8').Value = '  -0.40%  '
$ws.Range('E9').Value = '  +0.04%  '
$ws.Range('D10').Value = '''0.175'
$ws.Range('E10').Value = '  +2.41%  '
$ws.Range('E11').Value = '  +0.25%  '
$ws.Range('D12').Value = '''53.96'
$ws.Range('E12').Value = '  -1.93%  '
$ws.Range('E13').Value = '  -0.34%  '
$ws.Range('D14').Value = '''9.25'
$ws.Range('E14').Value = '  +0.95%  '
$ws.Range('D15').Value = '''3.928.02'
$ws.Range('E15').Value = '  -0.50%  '
$ws.Range('E16').Value = '  -1.01%  '
$ws.Range('E17').Value = '  +0.40%  '
$ws.Range('D18').Value = '''3.383.49'
$ws.Range('E18').Value = '  -0.92%  '
$ws.Range('D19').Value = '''65.770.74'
$ws.Range('E19').Value = '  +0.81%  '
$ws.Range('D20').Value = '''11.91'
$ws.Range('E20').Value = '  +0.06%  '
$ws.Range('D21').Value = '''0.998'
$ws.Range('E21').Value = '  +0.16%  '
$ws.Range('D22').Value = '''462.99'
$ws.Range('E22').Value = '  -2.25%  '
$ws.Range('D23').Value = '''4.94'
$ws.Range('E23').Value = '  -0.89%  '
$ws.Range('D24').Value = '''14.63'
$ws.Range('E24').Value = '  +8.31%  '
$ws.Range('D25').Value = '''89.42'
$ws.Range('E25').Value = '  +2.47%  '
$ws.Range('E26').Value = '  -1.21%  '
$ws.Range('D27').Value = '''2.93'
$ws.Range('E27').Value = '  -0.15%  '
$ws.Range('D28').Value = '''10.66'
$ws.Range('E28').Value = '  -2.52%  '
$ws.Range('D29').Value = '''8.74'
$ws.Range('E29').Value = '  -1.17%  '
$ws.Range('E30').Value = '  -0.61%  '
$ws.Range('D31').Value = '''6.62'
$ws.Range('E31').Value = '  -1.67%  '
$ws.Range('D32').Value = '''11.49'
$ws.Range('E32').Value = '  -0.74%  '
$ws.Range('D33').Value = '''581.30'
$ws.Range('E33').Value = '  +0.99%  '
$ws.Range('D34').Value = '''62.34'
$ws.Range('E34').Value = '  +0.37%  '
$ws.Range('E35').Value = '  -0.37%  '
$ws.Range('E36').Value = '  +0.08%  '
$ws.Range('D37').Value = '''3.59'
$ws.Range('E37').Value = '  +1.56%  '
$ws.Range('E38').Value = '  +1.38%  '
$ws.Range('D39').Value = '''36.01'
$ws.Range('E40').Value = '  +1.19%  '
$ws.Range('D41').Value = '''0.0₃0748'
$ws.Range('E41').Value = '  -2.01%  '
$ws.Range('D42').Value = '''3.107.23'
$ws.Range('E42').Value = '  +0.19%  '
$ws.Range('E43').Value = '  -1.00%  '
$ws.Range('E44').Value = '  +0.13%  '
$ws.Range('E45').Value = '  -1.17%  '
$ws.Range('D46').Value = '''2.45'
$ws.Range('E46').Value = '  -1.57%  '
$ws.Range('E47').Value = '  +0.08%  '
$ws.Range('D48').Value = '''1.00'
$ws.Range('E48').Value = '  +0.09%  '
$ws.Range('D49').Value = '''141.12'
$ws.Range('E49').Value = '  +2.51%  '
$ws.Range('B50').Value = 'dogwifhat'
$ws.Range('C50').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D50').Value = '''2.58'
$ws.Range('E50').Value = '  +9.33%  '
$ws.Range('B51').Value = 'THORChain'
$ws.Range('C51').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D51').Value = '''8.52'
$ws.Range('E51').Value = '  +2.10%  '
